$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '92.014.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.63%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.139.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.36%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.69%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '624.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.25%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.17'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.33%  '

# Row 8
$ws.Range("E8").Value = '  +2.32%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.05%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.128.22'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.05%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.755'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.98%  '

# Row 12
$ws.Range("E12").Value = '  +4.45%  '

# Row 13
$ws.Range("E13").Value = '  -1.33%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.56'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.97%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.08%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.121.55'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.16%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.688.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.98%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.088.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.25%  '

# Row 19
$ws.Range("E19").Value = '  -1.45%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.54%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000213'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.13%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.22%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '445.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.75%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.20%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.97%  '

# Row 26
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '90.56'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.23%  '

# Row 27
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.62%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.246.25'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.23%  '

# Row 29
$ws.Range("E29").Value = '  -0.10%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.251'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +28.74%  '

# Row 31
$ws.Range("E31").Value = '  +9.74%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.123'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +40.46%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.32'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.28%  '

# Row 34
$ws.Range("E34").Value = '  +11.01%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +12.93%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +27.69%  '

# Row 37
$ws.Range("E37").Value = '  +3.51%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '498.66'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.78%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.11%  '

# Row 40
$ws.Range("E40").Value = '  +1.50%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.26%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.426'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.83%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.24'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.26%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.95'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.13%  '

# Row 46
$ws.Range("E46").Value = '  +2.39%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '154.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.68%  '

# Row 48
$ws.Range("B48").Value = 'Binance-PegBSC-USD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.635'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -28.87%  '

# Row 49
$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.15%  '

# Row 50
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.44'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.32%  '

# Row 51
$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.37'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.24%  '
